$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.229.26"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.923.91"
$ws.Range("E3").Value = "  +3.65%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.42"
$ws.Range("E5").Value = "  +8.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.82"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.923.46"
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.432"
$ws.Range("E11").Value = "  +16.70%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.464.12"
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.134.90"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.03"
$ws.Range("E16").Value = "  +4.43%  "
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.931.71"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  +5.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.73"
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.79"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("E23").Value = "  +5.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.54"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.053.76"
$ws.Range("E26").Value = "  +3.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.24"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.73"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +3.89%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "507.22"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.76"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.09"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.24"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.63"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  +23.79%  "
$ws.Range("E40").Value = "  -4.43%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "182.17"
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +5.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.02"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.67"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.01"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.574"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.72"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.657"
$ws.Range("E51").Value = "  +3.17%  "
